# Swap the presentation's theme colour scheme ("Integral") for the
# generic Office default colour scheme ("Office").
#
# ppt/theme/theme1.xml (the theme used by the one slide master / all
# slides) only differs from ppt/theme/theme2.xml (the theme used by the
# notes master) in the <a:clrScheme> colour values - the font scheme and
# format scheme are byte-for-byte identical between the two themes, so
# the only substantive change is the 12 theme colours.
#
# PowerPoint's object model exposes those 12 colours as
# Slide.ThemeColorScheme.Item(1..12).RGB, in the fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink  12 folHlink
# RGB values are packed as 0xBBGGRR (PowerPoint's native colour order).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000  # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444  # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED  # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244  # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70  # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305  # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95  # folHlink -> 954F72
